$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: a batch of matches got re-scraped with swapped ordering (the two
# fixtures played on the same matchday ended up listed in the opposite row
# order vs. before). Columns F:V (home team .. url) swap between the two
# rows of each pair; columns A:E (Indice, pais, torneio, temporada,
# data_partida) stay anchored to their row.
# ---------------------------------------------------------------------------

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Swap-Rows($r1, $r2) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

Swap-Rows 22 23
Swap-Rows 42 43
Swap-Rows 60 61
Swap-Rows 78 79
Swap-Rows 84 85
Swap-Rows 101 102
Swap-Rows 107 108
Swap-Rows 127 128
Swap-Rows 136 137

# ---------------------------------------------------------------------------
# Part 2: two newly played fixtures got appended to the bottom of the sheet.
# Copy formatting from the last existing data row (157) so the new rows'
# A (bold/border/center) and E (datetime) styling matches the rest of the
# table, then fill in the values.
# ---------------------------------------------------------------------------

$ws.Range("A157:E157").Copy() | Out-Null
$ws.Range("A158:E158").PasteSpecial(-4122) | Out-Null
$ws.Range("A159:E159").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

function Set-Row($r, $values) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $values[$i]
    }
}

Set-Row 158 @(
    157,
    "indonesia",
    "liga-1",
    "2023-2024",
    45233.375,
    "PSS Sleman",
    0,
    "Bali United",
    1,
    3.37,
    "01/11/2023 21:12",
    2.71,
    "03/11/2023 08:58",
    3.51,
    "01/11/2023 21:12",
    3.72,
    "03/11/2023 08:58",
    1.93,
    "01/11/2023 21:12",
    2.38,
    "03/11/2023 08:58",
    "https://www.betexplorer.com/football/indonesia/liga-1/pss-sleman-bali-united/MkDDCmc6/"
)

Set-Row 159 @(
    158,
    "indonesia",
    "liga-1",
    "2023-2024",
    45233.54166666666,
    "PSM Makassar",
    2,
    "Persija Jakarta",
    3,
    2.45,
    "02/11/2023 01:12",
    2.17,
    "03/11/2023 12:56",
    3.07,
    "02/11/2023 01:12",
    3.26,
    "03/11/2023 12:55",
    2.73,
    "02/11/2023 01:12",
    3.45,
    "03/11/2023 12:56",
    "https://www.betexplorer.com/football/indonesia/liga-1/psm-makassar-persija-jakarta/r51MARRI/"
)
